$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..53 (columns A,B,C,D,E).
# Row 1 (headers) is unchanged. A new row was inserted at the top of the
# data block (old row 2 -> new row 3, etc.) and a new row was appended at
# the bottom (new row 53). All rows' B/C/D/E values were recomputed.
$data = @(
    @(39400,2007,11.13090654781821,2008,$null),
    @(39583,2008,$null,2009,$null),
    @(39765,2008,4.672550446571067,2009,$null),
    @(39948,2009,$null,2010,$null),
    @(40130,2009,-14.45332333832743,2010,$null),
    @(40310,2010,4.317520552235576,2011,3.856664261949816),
    @(40494,2010,8.600536527919633,2011,9.002271992040312),
    @(40676,2011,10.48563750975209,2012,5.658864198748459),
    @(40862,2011,10.25770250047622,2012,7.550992341868956),
    @(41044,2012,4.65880603412161,2013,5.409232631930561),
    @(41228,2012,4.639893381363169,2013,5.799303245920884),
    @(41409,2013,-0.8693696108860949,2014,2.311254283099773),
    @(41592,2013,0.3058963467304165,2014,2.638010271840896),
    @(41774,2014,4.057580120451165,2015,3.242274909585041),
    @(41957,2014,4.068173739091874,2015,5.156937396706884),
    @(42137,2015,3.942841799197594,2016,3.826653192455631),
    @(42321,2015,4.984288257750213,2016,3.4064284328156),
    @(42503,2016,2.435545128806416,2017,4.047617271894799),
    @(42689,2016,1.878184267712912,2017,2.129835064860464),
    @(42867,2017,4.218672000695523,2018,4.02808863798465),
    @(43053,2017,4.695933104194339,2018,5.022591279638045),
    @(43145,2018,7.356933926419673,2019,4.456100749631386),
    @(43235,2018,4.174017263680696,2019,2.074691389445271),
    @(43326,2018,5.32644934790627,2019,3.123778338720062),
    @(43418,2018,4.892602738886098,2019,0.6513682883433347),
    @(43510,2019,1.443764323860086,2020,2.608754752952369),
    @(43600,2019,1.699348375745302,2020,3.362240252406901),
    @(43691,2019,0.4648465628725118,2020,0.4030385202844711),
    @(43783,2019,0.8049382522247184,2020,3.036929265763488),
    @(43875,2020,0.642635895824295,2021,2.24075411271687),
    @(43966,2020,-2.856524424985296,2021,0.1239622353166103),
    @(44068,2020,-9.810777096850787,2021,-4.291312636561706),
    @(44159,2020,-8.784173899737169,2021,1.573231731123359),
    @(44251,2021,5.992199201897175,2022,3.59700023027214),
    @(44341,2021,6.240787792289715,2022,4.112367048326182),
    @(44432,2021,5.797134106720514,2022,4.062819790119954),
    @(44525,2021,5.110501195359984,2022,1.88131082127776),
    @(44617,2022,5.158753599182209,2023,3.207070318919869),
    @(44706,2022,3.815916106066686,2023,2.762663830671319),
    @(44798,2022,4.232564748995715,2023,2.257990710759383),
    @(44890,2022,5.120680133083599,2023,2.522545412785848),
    @(44981,2023,1.135459337362521,2024,2.746436727258894),
    @(45071,2023,0.757583445265464,2024,2.564429185896056),
    @(45163,2023,0.08070151925247959,2024,1.796717073915977),
    @(45254,2023,-0.5532735011319234,2024,-0.657715646732393),
    @(45345,2024,-2.103425609777143,2025,2.179038063425076),
    @(45436,2024,-0.152046383567539,2025,2.659925441240518),
    @(45534,2024,-0.9685570952743805,2025,1.778332450996523),
    @(45618,2024,-1.069674659641462,2025,0.5636794832278413),
    @(45713,2025,-2.92939668194816,2026,1.921265997709742),
    @(45800,2025,-1.93082584212636,2026,2.2373336846083),
    @(45891,2025,-2.436529450546909,2026,1.970944416887122)
)

# Row 53 is brand new territory on the sheet; give it the same format
# (date number format, border, font, alignment) as the existing date
# column cells before writing values, so the style matches the rest of
# column A (s="2" in the original date column style).
$ws.Cells.Item(52, 1).Copy() | Out-Null
$ws.Cells.Item(53, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}
